# Update valid username/password on the "Valid_Login" sheet and make it the active sheet.

$wb = $excel.ActiveWorkbook

$wsValid = $wb.Worksheets.Item("Valid_Login")

# Update the valid credentials (both the username and password cell now hold
# the same new value, per the source data change).
$wsValid.Range("A2").Value = "RajGuru11991145@gmail.com"
$wsValid.Range("B2").Value = "RajGuru11991145@gmail.com"

# Make Valid_Login the active sheet/tab and select cell A2 on it.
$wsValid.Activate()
$wsValid.Range("A2").Select()

$wb.Save()
